$d = $word.ActiveDocument

$d.Content.Find.Execute("2025-02-02 Sunday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-02-03 Monday", 2) | Out-Null
$d.Content.Find.Execute("75+8=83", $true, $false, $false, $false, $false, $true, 1, $false, "34+42=76", 2) | Out-Null
$d.Content.Find.Execute("57+26=83", $true, $false, $false, $false, $false, $true, 1, $false, "32+19=51", 2) | Out-Null
$d.Content.Find.Execute("26+54=80", $true, $false, $false, $false, $false, $true, 1, $false, "94-24=70", 2) | Out-Null
$d.Content.Find.Execute("29-22=7", $true, $false, $false, $false, $false, $true, 1, $false, "26+55=81", 2) | Out-Null
$d.Content.Find.Execute("87-64=23", $true, $false, $false, $false, $false, $true, 1, $false, "72+7=79", 2) | Out-Null
$d.Content.Find.Execute("21+36=57", $true, $false, $false, $false, $false, $true, 1, $false, "1+86=87", 2) | Out-Null
$d.Content.Find.Execute("68-6=62", $true, $false, $false, $false, $false, $true, 1, $false, "0+69=69", 2) | Out-Null
$d.Content.Find.Execute("60+19=79", $true, $false, $false, $false, $false, $true, 1, $false, "58-52=6", 2) | Out-Null
$d.Content.Find.Execute("42+25=67", $true, $false, $false, $false, $false, $true, 1, $false, "30+63=93", 2) | Out-Null
$d.Content.Find.Execute("19+34=53", $true, $false, $false, $false, $false, $true, 1, $false, "97-28=69", 2) | Out-Null
$d.Content.Find.Execute("95-75=20", $true, $false, $false, $false, $false, $true, 1, $false, "51-4=47", 2) | Out-Null
$d.Content.Find.Execute("0+61=61", $true, $false, $false, $false, $false, $true, 1, $false, "55-8=47", 2) | Out-Null
$d.Content.Find.Execute("3+51=54", $true, $false, $false, $false, $false, $true, 1, $false, "28+20=48", 2) | Out-Null
$d.Content.Find.Execute("5+87=92", $true, $false, $false, $false, $false, $true, 1, $false, "11+53=64", 2) | Out-Null
$d.Content.Find.Execute("98-89=9", $true, $false, $false, $false, $false, $true, 1, $false, "65+5=70", 2) | Out-Null
$d.Content.Find.Execute("1+45=46", $true, $false, $false, $false, $false, $true, 1, $false, "94-48=46", 2) | Out-Null
$d.Content.Find.Execute("24+11=35", $true, $false, $false, $false, $false, $true, 1, $false, "23+12=35", 2) | Out-Null
$d.Content.Find.Execute("33+62=95", $true, $false, $false, $false, $false, $true, 1, $false, "63+32=95", 2) | Out-Null
$d.Content.Find.Execute("7-7=0", $true, $false, $false, $false, $false, $true, 1, $false, "94-5=89", 2) | Out-Null
$d.Content.Find.Execute("2+24=26", $true, $false, $false, $false, $false, $true, 1, $false, "34+23=57", 2) | Out-Null
$d.Content.Find.Execute("33-13=20", $true, $false, $false, $false, $false, $true, 1, $false, "7+89=96", 2) | Out-Null
$d.Content.Find.Execute("40+38=78", $true, $false, $false, $false, $false, $true, 1, $false, "43+19=62", 2) | Out-Null
$d.Content.Find.Execute("78-35=43", $true, $false, $false, $false, $false, $true, 1, $false, "50+13=63", 2) | Out-Null
$d.Content.Find.Execute("91-72=19", $true, $false, $false, $false, $false, $true, 1, $false, "80-18=62", 2) | Out-Null
$d.Content.Find.Execute("93-55=38", $true, $false, $false, $false, $false, $true, 1, $false, "77-68=9", 2) | Out-Null
$d.Content.Find.Execute("88+7=95", $true, $false, $false, $false, $false, $true, 1, $false, "26+8=34", 2) | Out-Null
$d.Content.Find.Execute("73-8=65", $true, $false, $false, $false, $false, $true, 1, $false, "19-18=1", 2) | Out-Null
$d.Content.Find.Execute("54+39=93", $true, $false, $false, $false, $false, $true, 1, $false, "21+45=66", 2) | Out-Null
$d.Content.Find.Execute("10+37=47", $true, $false, $false, $false, $false, $true, 1, $false, "59-18=41", 2) | Out-Null
$d.Content.Find.Execute("29-13=16", $true, $false, $false, $false, $false, $true, 1, $false, "18+81=99", 2) | Out-Null
$d.Content.Find.Execute("43-4=39", $true, $false, $false, $false, $false, $true, 1, $false, "15+13=28", 2) | Out-Null
$d.Content.Find.Execute("70-62=8", $true, $false, $false, $false, $false, $true, 1, $false, "80+11=91", 2) | Out-Null
$d.Content.Find.Execute("71+26=97", $true, $false, $false, $false, $false, $true, 1, $false, "44-9=35", 2) | Out-Null
$d.Content.Find.Execute("98-10=88", $true, $false, $false, $false, $false, $true, 1, $false, "26+68=94", 2) | Out-Null
$d.Content.Find.Execute("76-8=68", $true, $false, $false, $false, $false, $true, 1, $false, "93-45=48", 2) | Out-Null
$d.Content.Find.Execute("61-47=14", $true, $false, $false, $false, $false, $true, 1, $false, "38+18=56", 2) | Out-Null
$d.Content.Find.Execute("38+2=40", $true, $false, $false, $false, $false, $true, 1, $false, "59+14=73", 2) | Out-Null
$d.Content.Find.Execute("89-46=43", $true, $false, $false, $false, $false, $true, 1, $false, "41+34=75", 2) | Out-Null
$d.Content.Find.Execute("74-46=28", $true, $false, $false, $false, $false, $true, 1, $false, "82-52=30", 2) | Out-Null
$d.Content.Find.Execute("33-10=23", $true, $false, $false, $false, $false, $true, 1, $false, "96+0=96", 2) | Out-Null
$d.Content.Find.Execute("52-30=22", $true, $false, $false, $false, $false, $true, 1, $false, "19+42=61", 2) | Out-Null
$d.Content.Find.Execute("99-93=6", $true, $false, $false, $false, $false, $true, 1, $false, "56+32=88", 2) | Out-Null
$d.Content.Find.Execute("64-3=61", $true, $false, $false, $false, $false, $true, 1, $false, "36+3=39", 2) | Out-Null
$d.Content.Find.Execute("69-66=3", $true, $false, $false, $false, $false, $true, 1, $false, "42-29=13", 2) | Out-Null
$d.Content.Find.Execute("23+9=32", $true, $false, $false, $false, $false, $true, 1, $false, "90-36=54", 2) | Out-Null
$d.Content.Find.Execute("97-17=80", $true, $false, $false, $false, $false, $true, 1, $false, "27+45=72", 2) | Out-Null
$d.Content.Find.Execute("47+9=56", $true, $false, $false, $false, $false, $true, 1, $false, "95-1=94", 2) | Out-Null
$d.Content.Find.Execute("69+1=70", $true, $false, $false, $false, $false, $true, 1, $false, "2+59=61", 2) | Out-Null
$d.Content.Find.Execute("3+68=71", $true, $false, $false, $false, $false, $true, 1, $false, "26+55=81", 2) | Out-Null
$d.Content.Find.Execute("30+41=71", $true, $false, $false, $false, $false, $true, 1, $false, "32-9=23", 2) | Out-Null
$d.Content.Find.Execute("74-19=55", $true, $false, $false, $false, $false, $true, 1, $false, "42+36=78", 2) | Out-Null
$d.Content.Find.Execute("50+36=86", $true, $false, $false, $false, $false, $true, 1, $false, "80-53=27", 2) | Out-Null
$d.Content.Find.Execute("26+7=33", $true, $false, $false, $false, $false, $true, 1, $false, "44-7=37", 2) | Out-Null
$d.Content.Find.Execute("23+76=99", $true, $false, $false, $false, $false, $true, 1, $false, "91-34=57", 2) | Out-Null
$d.Content.Find.Execute("85-18=67", $true, $false, $false, $false, $false, $true, 1, $false, "92-63=29", 2) | Out-Null
$d.Content.Find.Execute("93+6=99", $true, $false, $false, $false, $false, $true, 1, $false, "10+46=56", 2) | Out-Null
$d.Content.Find.Execute("91-87=4", $true, $false, $false, $false, $false, $true, 1, $false, "32+46=78", 2) | Out-Null
$d.Content.Find.Execute("18+66=84", $true, $false, $false, $false, $false, $true, 1, $false, "6+73=79", 2) | Out-Null
$d.Content.Find.Execute("72-66=6", $true, $false, $false, $false, $false, $true, 1, $false, "37+59=96", 2) | Out-Null
$d.Content.Find.Execute("62-51=11", $true, $false, $false, $false, $false, $true, 1, $false, "42-21=21", 2) | Out-Null
$d.Content.Find.Execute("62-13=49", $true, $false, $false, $false, $false, $true, 1, $false, "18+77=95", 2) | Out-Null
$d.Content.Find.Execute("37+9=46", $true, $false, $false, $false, $false, $true, 1, $false, "71-22=49", 2) | Out-Null
$d.Content.Find.Execute("83+13=96", $true, $false, $false, $false, $false, $true, 1, $false, "89-11=78", 2) | Out-Null
$d.Content.Find.Execute("16+49=65", $true, $false, $false, $false, $false, $true, 1, $false, "79-66=13", 2) | Out-Null
$d.Content.Find.Execute("30-14=16", $true, $false, $false, $false, $false, $true, 1, $false, "86-10=76", 2) | Out-Null
$d.Content.Find.Execute("19+45=64", $true, $false, $false, $false, $false, $true, 1, $false, "17+10=27", 2) | Out-Null
$d.Content.Find.Execute("82-55=27", $true, $false, $false, $false, $false, $true, 1, $false, "9+31=40", 2) | Out-Null
$d.Content.Find.Execute("5+3=8", $true, $false, $false, $false, $false, $true, 1, $false, "84-17=67", 2) | Out-Null
$d.Content.Find.Execute("73+22=95", $true, $false, $false, $false, $false, $true, 1, $false, "76+19=95", 2) | Out-Null
$d.Content.Find.Execute("90-11=79", $true, $false, $false, $false, $false, $true, 1, $false, "68-65=3", 2) | Out-Null
$d.Content.Find.Execute("50-16=34", $true, $false, $false, $false, $false, $true, 1, $false, "10+59=69", 2) | Out-Null
$d.Content.Find.Execute("0+9=9", $true, $false, $false, $false, $false, $true, 1, $false, "71-10=61", 2) | Out-Null
$d.Content.Find.Execute("25+69=94", $true, $false, $false, $false, $false, $true, 1, $false, "44+28=72", 2) | Out-Null
$d.Content.Find.Execute("3+88=91", $true, $false, $false, $false, $false, $true, 1, $false, "69-8=61", 2) | Out-Null
$d.Content.Find.Execute("30+28=58", $true, $false, $false, $false, $false, $true, 1, $false, "61-8=53", 2) | Out-Null
$d.Content.Find.Execute("41+18=59", $true, $false, $false, $false, $false, $true, 1, $false, "84-56=28", 2) | Out-Null
$d.Content.Find.Execute("25+16=41", $true, $false, $false, $false, $false, $true, 1, $false, "75-57=18", 2) | Out-Null
$d.Content.Find.Execute("28+19=47", $true, $false, $false, $false, $false, $true, 1, $false, "0+14=14", 2) | Out-Null
$d.Content.Find.Execute("1+60=61", $true, $false, $false, $false, $false, $true, 1, $false, "21+66=87", 2) | Out-Null
$d.Content.Find.Execute("51+24=75", $true, $false, $false, $false, $false, $true, 1, $false, "94-6=88", 2) | Out-Null
$d.Content.Find.Execute("21+51=72", $true, $false, $false, $false, $false, $true, 1, $false, "31+17=48", 2) | Out-Null
$d.Content.Find.Execute("78-65=13", $true, $false, $false, $false, $false, $true, 1, $false, "74-24=50", 2) | Out-Null
$d.Content.Find.Execute("84-37=47", $true, $false, $false, $false, $false, $true, 1, $false, "14+78=92", 2) | Out-Null
$d.Content.Find.Execute("14+5=19", $true, $false, $false, $false, $false, $true, 1, $false, "38+7=45", 2) | Out-Null
$d.Content.Find.Execute("34+31=65", $true, $false, $false, $false, $false, $true, 1, $false, "25+50=75", 2) | Out-Null
$d.Content.Find.Execute("45+20=65", $true, $false, $false, $false, $false, $true, 1, $false, "83-63=20", 2) | Out-Null
$d.Content.Find.Execute("8+25=33", $true, $false, $false, $false, $false, $true, 1, $false, "17+45=62", 2) | Out-Null
$d.Content.Find.Execute("87+0=87", $true, $false, $false, $false, $false, $true, 1, $false, "50-24=26", 2) | Out-Null
$d.Content.Find.Execute("85-10=75", $true, $false, $false, $false, $false, $true, 1, $false, "97-36=61", 2) | Out-Null
$d.Content.Find.Execute("2+86=88", $true, $false, $false, $false, $false, $true, 1, $false, "99-33=66", 2) | Out-Null
$d.Content.Find.Execute("70-9=61", $true, $false, $false, $false, $false, $true, 1, $false, "69+11=80", 2) | Out-Null
$d.Content.Find.Execute("87-75=12", $true, $false, $false, $false, $false, $true, 1, $false, "8+58=66", 2) | Out-Null
$d.Content.Find.Execute("69-3=66", $true, $false, $false, $false, $false, $true, 1, $false, "11+48=59", 2) | Out-Null
$d.Content.Find.Execute("17+34=51", $true, $false, $false, $false, $false, $true, 1, $false, "73-25=48", 2) | Out-Null
$d.Content.Find.Execute("28+22=50", $true, $false, $false, $false, $false, $true, 1, $false, "0+30=30", 2) | Out-Null
$d.Content.Find.Execute("27+12=39", $true, $false, $false, $false, $false, $true, 1, $false, "71-70=1", 2) | Out-Null
$d.Content.Find.Execute("86-12=74", $true, $false, $false, $false, $false, $true, 1, $false, "5+66=71", 2) | Out-Null
$d.Content.Find.Execute("34+61=95", $true, $false, $false, $false, $false, $true, 1, $false, "78-3=75", 2) | Out-Null
$d.Content.Find.Execute("69-35=34", $true, $false, $false, $false, $false, $true, 1, $false, "70-4=66", 2) | Out-Null
$d.Content.Find.Execute("5+2=7", $true, $false, $false, $false, $false, $true, 1, $false, "68-56=12", 2) | Out-Null
